# Auto-generated-assisted edit script for 北京-漫展信息.xlsx
# Applies: sheet1 (展览) F-column bumps + 2 new rows inserted at 34 and 39 (rows shift down),
# sheet2 (演出), sheet3 (本地生活), sheet4 (全部类型) F/I-column updates.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 展览

# --- F-column ('想去人数') bumps for existing rows 4-33 ---
$sheet1FChanges = @{
    4 = 552
    6 = 1593
    9 = 720
    10 = 2662
    11 = 2662
    12 = 14
    13 = 1732
    14 = 603
    15 = 266
    16 = 682
    17 = 4955
    18 = 164
    19 = 62
    20 = 689
    21 = 3381
    22 = 854
    23 = 30
    24 = 67
    26 = 2405
    28 = 362
    29 = 19
    32 = 1287
    33 = 799
}
foreach ($r in $sheet1FChanges.Keys) {
    $ws1.Cells.Item([int]$r, 6).Value = $sheet1FChanges[$r]
}

# --- Rows 34-41: two new rows inserted (at original position 34 and 39); ---
# --- remaining old rows shift down and carry small F-column bumps.        ---
# --- Write the fully resolved final content for rows 34..41 directly.     ---
$sheet1Rows = @{}
$sheet1Rows[34] = @{
    A = 33
    B = '2024-12-07'
    C = '北京·CCW华彩国潮动漫游戏世界1.0'
    D = '北花园路1号超级蜂巢C座 超级蜂巢国际会议中心'
    E = '2024.12.07 09:00-12.08 17:00'
    F = 0
    G = 55
    H = 'https://show.bilibili.com/platform/detail.html?id=93741'
    I = '//i0.hdslb.com/bfs/openplatform/202410/n42wlerB1729262442496.png'
}
$sheet1Rows[35] = @{
    A = 34
    B = '2024-12-07'
    C = '北京·排球少年同人ONLY'
    D = '永外高庄138号 北京大红门国际会展中心'
    E = '2024.12.07 10:00-12.07 17:00'
    F = 54
    G = 60
    H = 'https://show.bilibili.com/platform/detail.html?id=93501'
    I = '//i1.hdslb.com/bfs/openplatform/202410/Nq2SuSWE1728971583727.jpeg'
}
$sheet1Rows[36] = @{
    A = 35
    B = '2024-12-14'
    C = '北京·thebONE×GOJO超次元动漫游戏嘉年华17th'
    D = 'B1层西区会员活动中心（阳坊涮肉对面） 北投购物公园'
    E = '2024.12.14 09:30-12.15 17:00'
    F = 17
    G = 1
    H = 'https://show.bilibili.com/platform/detail.html?id=93495'
    I = '//i0.hdslb.com/bfs/openplatform/202410/Fzz24Usj1728969298701.jpeg'
}
$sheet1Rows[37] = @{
    A = 36
    B = '2024-12-14'
    C = '北京·奇想派对第五届'
    D = '学清路38号金码大厦B座(六道口地铁站B东北口步行110米) BOM嘻番里'
    E = '2024.12.14 10:00-12.15 17:30'
    F = 53
    G = 45
    H = 'https://show.bilibili.com/platform/detail.html?id=91077'
    I = '//i1.hdslb.com/bfs/openplatform/202408/zMayUoC81724229782742.jpeg'
}
$sheet1Rows[38] = @{
    A = 37
    B = '2024-12-28'
    C = '北京·第20届IJOY漫展xCGF游戏节'
    D = '天辰东路7号 北京国家会议中心'
    E = '2024.12.28 09:00-12.29 17:00'
    F = 1404
    G = 8.800000000000001
    H = 'https://show.bilibili.com/platform/detail.html?id=92633'
    I = '//i0.hdslb.com/bfs/openplatform/202409/EQg8HwjJ1726734597607.jpeg'
}
$sheet1Rows[39] = @{
    A = 38
    B = '2025-01-01'
    C = '北京·第五人格only同人展2.0'
    D = '北花园路1号超级蜂巢C座 超级蜂巢国际会议中心'
    E = '2025.01.01 10:00-01.01 17:00'
    F = 1
    G = 68
    H = 'https://show.bilibili.com/platform/detail.html?id=93734'
    I = '//i2.hdslb.com/bfs/openplatform/202410/AwXBn8Jn1729336699687.jpeg'
}
$sheet1Rows[40] = @{
    A = 39
    B = '2025-01-17'
    C = ' 北京·第21届IJOY漫展xCGF游戏节'
    D = '天辰东路7号 北京国家会议中心'
    E = '2025.01.17 09:00-01.19 17:00'
    F = 1369
    G = 8.800000000000001
    H = 'https://show.bilibili.com/platform/detail.html?id=92634'
    I = '//i0.hdslb.com/bfs/openplatform/202409/ASXIizNW1726735204415.jpeg'
}
$sheet1Rows[41] = @{
    A = 40
    B = '2025-04-19'
    C = '北京·可行中国动漫游戏节'
    D = '焦化路甲18号 东进国际中心'
    E = '2025.04.19 09:00-04.20 18:00'
    F = 90
    G = 85
    H = 'https://show.bilibili.com/platform/detail.html?id=92495'
    I = '//i1.hdslb.com/bfs/openplatform/202409/28QBTqAo1726293348310.jpeg'
}

foreach ($r in 34..41) {
    $row = $sheet1Rows[$r]
    $ws1.Cells.Item($r, 1).Value = [double]$row['A']
    # Column B holds a literal date-like string (e.g. '2024-12-07'); prefix with an
    # apostrophe so Excel stores literal text instead of auto-converting to a date serial.
    $ws1.Cells.Item($r, 2).Value = "'" + $row['B']
    $ws1.Cells.Item($r, 3).Value = $row['C']
    $ws1.Cells.Item($r, 4).Value = $row['D']
    $ws1.Cells.Item($r, 5).Value = $row['E']
    $ws1.Cells.Item($r, 6).Value = [double]$row['F']
    $ws1.Cells.Item($r, 7).Value = [double]$row['G']
    $ws1.Cells.Item($r, 8).Value = $row['H']
    $ws1.Cells.Item($r, 9).Value = $row['I']
}

# Rows 40-41 are brand-new sheet rows (sheet used to end at row 39); column A there
# needs the same bold/centered/bordered style used by the rest of the index column.
# Copy number/cell formatting only (not values) from an existing index cell.
$ws1.Range("A33").Copy()
$ws1.Range("A40:A41").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)   # 演出
$ws2.Range('I15').Value = '//i0.hdslb.com/bfs/openplatform/202410/a4qEFGiD1729491938007.png'
$ws2.Range('F16').Value = 137
$ws2.Range('F18').Value = 257

# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)   # 本地生活
$ws3.Range('F3').Value = 842
$ws3.Range('F4').Value = 240
$ws3.Range('F6').Value = 12
$ws3.Range('I6').Value = '//i0.hdslb.com/bfs/openplatform/202410/tBaCv15Q1729493037977.png'
$ws3.Range('F7').Value = 12
$ws3.Range('I7').Value = '//i0.hdslb.com/bfs/openplatform/202410/lVX1X9gy1729492890064.png'
$ws3.Range('F8').Value = 2

# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)   # 全部类型
$ws4.Range('F6').Value = 842
$ws4.Range('F7').Value = 240
$ws4.Range('F10').Value = 552
$ws4.Range('F11').Value = 12
$ws4.Range('I11').Value = '//i0.hdslb.com/bfs/openplatform/202410/tBaCv15Q1729493037977.png'
$ws4.Range('F12').Value = 12
$ws4.Range('I12').Value = '//i0.hdslb.com/bfs/openplatform/202410/lVX1X9gy1729492890064.png'
$ws4.Range('F17').Value = 1593
$ws4.Range('F21').Value = 2662
$ws4.Range('F23').Value = 1732
$ws4.Range('F25').Value = 603
$ws4.Range('F26').Value = 266
$ws4.Range('F27').Value = 682
$ws4.Range('F28').Value = 4955
$ws4.Range('F29').Value = 62
$ws4.Range('F30').Value = 689
$ws4.Range('F31').Value = 3381
$ws4.Range('F32').Value = 854
$ws4.Range('F33').Value = 67
$ws4.Range('F36').Value = 2405
$ws4.Range('F38').Value = 362
$ws4.Range('F39').Value = 19
$ws4.Range('F42').Value = 1287
$ws4.Range('F43').Value = 137
$ws4.Range('F44').Value = 257
$ws4.Range('F46').Value = 799
$ws4.Range('F47').Value = 54
$ws4.Range('F50').Value = 1404

Write-Host 'edit.ps1 completed'
